$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Values are forced to text (leading apostrophe) so Excel doesn't reinterpret
# numeric-looking strings (e.g. "10.20", "1.00", "68.350.13") as numbers,
# then the cell style is reset to Normal to drop the quote-prefix flag Excel
# adds, keeping formatting identical to the original inline-string cells.

$ws.Range("D2").Value = "'68.350.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'3.766.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.35%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'595.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.22%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'168.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.47%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.766.81"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.21%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.42%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.69%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.40%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -2.40%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -3.71%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'36.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.11%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.399.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.50%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.763.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.66%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'68.336.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.80%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'18.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.47%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.85%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -0.28%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'467.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.37%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -3.22%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'84.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.23%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -4.79%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.09%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.55%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.39%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.14%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.914.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.55%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.43%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'7.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.58%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.49%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'30.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.34%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'9.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.80%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'3.721.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.85%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -3.42%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -10.47%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.50%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.56%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'5.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.48%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.08%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.01%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -2.93%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.52%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'43.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +11.50%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.60%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'406.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.15%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'45.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'145.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.15%  "
$ws.Range("E51").Style = "Normal"
